$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1498.1666
$ws.Range("I19").Value = 1363.6666
$ws.Range("J19").Value = 1632.6666
$ws.Range("K19").Value = 1363.6666
$ws.Range("L19").Value = 1632.6666
$ws.Range("M19").Value = -1188.6666
$ws.Range("N19").Value = -1982.6666
$ws.Range("H31").Value = 8333.666999999999
$ws.Range("I31").Value = 7001
$ws.Range("K31").Value = 21003
$ws.Range("M31").Value = -20773
$ws.Range("H69").Value = 21363.182
$ws.Range("I69").Value = 9995
$ws.Range("J69").Value = 22500
$ws.Range("K69").Value = 29985
$ws.Range("L69").Value = 67500
$ws.Range("M69").Value = -29111
$ws.Range("N69").Value = -69248
$ws.Range("H72").Value = 21363.182
$ws.Range("I72").Value = 9995
$ws.Range("J72").Value = 22500
$ws.Range("K72").Value = 89955
$ws.Range("L72").Value = 202500
$ws.Range("M72").Value = -85587
$ws.Range("N72").Value = -211236
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H92").Value = 918.4666999999999
$ws.Range("I92").Value = 841.2143
$ws.Range("J92").Value = 2000
$ws.Range("K92").Value = 841.2143
$ws.Range("L92").Value = 2000
$ws.Range("M92").Value = 406.7857
$ws.Range("N92").Value = -4496
$ws.Range("H125").Value = 3395110
$ws.Range("I125").Value = 4241788
$ws.Range("K125").Value = 38176092
$ws.Range("M125").Value = -38173632
$ws.Range("H129").Value = 8959.833000000001
$ws.Range("I129").Value = 1652.1
$ws.Range("K129").Value = 4956.299999999999
$ws.Range("M129").Value = 43.70000000000073
$ws.Range("H137").Value = 3041.0454
$ws.Range("I137").Value = 2576.7646
$ws.Range("K137").Value = 7730.293799999999
$ws.Range("M137").Value = -5180.293799999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2161.0637
$ws.Range("I32").Value = 2152.6047
$ws.Range("K32").Value = 2152.6047
$ws.Range("M32").Value = -1865.6047
$ws.Range("H74").Value = 34485744
$ws.Range("I74").Value = 37039504
$ws.Range("J74").Value = 10000
$ws.Range("K74").Value = 37039504
$ws.Range("L74").Value = 10000
$ws.Range("M74").Value = -37038630
$ws.Range("N74").Value = -11748
$ws.Range("H77").Value = 34485744
$ws.Range("I77").Value = 37039504
$ws.Range("J77").Value = 10000
$ws.Range("K77").Value = 185197520
$ws.Range("L77").Value = 50000
$ws.Range("M77").Value = -185193152
$ws.Range("N77").Value = -58736
$ws.Range("H132").Value = 5266771.5
$ws.Range("I132").Value = 6253460
$ws.Range("J132").Value = 4433.3335
$ws.Range("K132").Value = 18760380
$ws.Range("L132").Value = 13300.0005
$ws.Range("M132").Value = -18757850
$ws.Range("N132").Value = -18360.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1352
$ws.Range("I20").Value = 1469.3334
$ws.Range("J20").Value = 1000
$ws.Range("K20").Value = 1469.3334
$ws.Range("L20").Value = 1000
$ws.Range("M20").Value = -1222.3334
$ws.Range("N20").Value = -1494
$ws.Range("H80").Value = 55556570
$ws.Range("I80").Value = 1299
$ws.Range("J80").Value = 62500976
$ws.Range("K80").Value = 1299
$ws.Range("L80").Value = 62500976
$ws.Range("M80").Value = -301
$ws.Range("N80").Value = -62502972
$ws.Range("H83").Value = 55556570
$ws.Range("I83").Value = 1299
$ws.Range("J83").Value = 62500976
$ws.Range("K83").Value = 6495
$ws.Range("L83").Value = 312504880
$ws.Range("M83").Value = -1503
$ws.Range("N83").Value = -312514864

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6338.4
$ws.Range("I31").Value = 4159.6924
$ws.Range("K31").Value = 4159.6924
$ws.Range("M31").Value = -3864.6924
$ws.Range("H34").Value = 6338.4
$ws.Range("I34").Value = 4159.6924
$ws.Range("K34").Value = 4159.6924
$ws.Range("M34").Value = -3957.6924
$ws.Range("H86").Value = 9812.895
$ws.Range("I86").Value = 7113.273
$ws.Range("K86").Value = 7113.273
$ws.Range("M86").Value = -5990.273
$ws.Range("H89").Value = 9812.895
$ws.Range("I89").Value = 7113.273
$ws.Range("K89").Value = 35566.365
$ws.Range("M89").Value = -29950.365
$ws.Range("H105").Value = 5630958
$ws.Range("I105").Value = 5630958
$ws.Range("K105").Value = 5630958
$ws.Range("M105").Value = -5629211
$ws.Range("H132").Value = 200003200
$ws.Range("I132").Value = 250002990
$ws.Range("J132").Value = 4014
$ws.Range("K132").Value = 750008970
$ws.Range("L132").Value = 12042
$ws.Range("M132").Value = -750006440
$ws.Range("N132").Value = -17102

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 9999999
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H39").Value = 1289.6
$ws.Range("J39").Value = 1733.1666
$ws.Range("L39").Value = 5199.4998
$ws.Range("N39").Value = -5787.4998
$ws.Range("H93").Value = 9599.9
$ws.Range("J93").Value = 9599.9
$ws.Range("L93").Value = 28799.7
$ws.Range("N93").Value = -32543.7
$ws.Range("H94").Value = 19167.9
$ws.Range("J94").Value = 20742.111
$ws.Range("L94").Value = 62226.333
$ws.Range("N94").Value = -63578.333
$ws.Range("H102").Value = 2499
$ws.Range("I102").Value = 2499
$ws.Range("K102").Value = 7497
$ws.Range("M102").Value = -5063
$ws.Range("H107").Value = 1124.2
$ws.Range("I107").Value = 468.75
$ws.Range("J107").Value = 1873.2858
$ws.Range("K107").Value = 1406.25
$ws.Range("L107").Value = 5619.857400000001
$ws.Range("M107").Value = 513.75
$ws.Range("N107").Value = -9459.857400000001
$ws.Range("H140").Value = 2646.375
$ws.Range("I140").Value = 2595.8572
$ws.Range("J140").Value = 3000
$ws.Range("K140").Value = 7787.571599999999
$ws.Range("L140").Value = 9000
$ws.Range("M140").Value = -2607.571599999999
$ws.Range("N140").Value = -19360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3025
$ws.Range("I80").Value = 2644.5557
$ws.Range("K80").Value = 2644.5557
$ws.Range("M80").Value = -1646.5557
$ws.Range("H83").Value = 3025
$ws.Range("I83").Value = 2644.5557
$ws.Range("K83").Value = 13222.7785
$ws.Range("M83").Value = -8230.7785
$ws.Range("H92").Value = 2443.6667
$ws.Range("J92").Value = 2443.6667
$ws.Range("L92").Value = 2443.6667
$ws.Range("N92").Value = -6187.6667
$ws.Range("H111").Value = 38330
$ws.Range("J111").Value = 38330
$ws.Range("L111").Value = 38330
$ws.Range("N111").Value = -44464
$ws.Range("H113").Value = 59570.5
$ws.Range("I113").Value = 63009.94
$ws.Range("J113").Value = 1100
$ws.Range("K113").Value = 63009.94
$ws.Range("L113").Value = 1100
$ws.Range("M113").Value = -60839.94
$ws.Range("N113").Value = -5440
$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530
$ws.Range("H132").Value = 5686966
$ws.Range("I132").Value = 6584310.5
$ws.Range("K132").Value = 19752931.5
$ws.Range("M132").Value = -19750401.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H100").Value = 8679001
$ws.Range("I100").Value = 9073433
$ws.Range("K100").Value = 9073433
$ws.Range("M100").Value = -9072892
$ws.Range("H132").Value = 16561540
$ws.Range("I132").Value = 20879820
$ws.Range("J132").Value = 8136
$ws.Range("K132").Value = 62639460
$ws.Range("L132").Value = 24408
$ws.Range("M132").Value = -62636930
$ws.Range("N132").Value = -29468

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 5224995.5
$ws.Range("I3").Value = 5224995.5
$ws.Range("K3").Value = 5224995.5
$ws.Range("M3").Value = -5224881.5
$ws.Range("H8").Value = 295000000
$ws.Range("I8").Value = 295000000
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 295000000
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -294999860
$ws.Range("N8").ClearContents()
$ws.Range("H43").Value = 20030
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H46").Value = 50000
$ws.Range("J46").Value = 50000
$ws.Range("L46").Value = 50000
$ws.Range("N46").Value = -50462
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 150000
$ws.Range("N134").Value = -155070
$ws.Range("H136").Value = 27780438
$ws.Range("I136").Value = 31252756
$ws.Range("K136").Value = 93758268
$ws.Range("M136").Value = -93755718
